$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in A1 from "TestScenario" to "TestCases"
$ws.Range("A1").Value = "TestCases"

# Reset selection back to the default cell (A1) instead of B2
$ws.Range("A1").Select()
